$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.194.35'
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').Value = '2.055.54'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.616'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.67'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.80%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0828'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.51%  '
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.03%  '
$ws.Range('D13').Value = '2.359.70'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.761'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('D17').Value = '2.055.57'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '38.134.20'
$ws.Range('E18').Value = '  +2.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = '0.0₃0832'
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '225.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('E30').Value = '  -2.03%  '
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.29'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +11.04%  '
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('E38').Value = '  +0.75%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').Value = '1.526.72'
$ws.Range('E40').Value = '  +3.87%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0219'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.37%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.03'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.76'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.80%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0929'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('E47').Value = '  -6.37%  '
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.08'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '2.249.65'
$ws.Range('E51').Value = '  +1.34%  '
